$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B3: ImagePath -> ImageId ---
$ws.Range("B3").Value = "Dateisystem für Bilder, Product::ImageId, ProductForm & ProductOverview: Bild anzeigen"

# --- New column G header "Sonstiges" (bold black, centered) ---
$ws.Range("G1").Value = "Sonstiges"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108

# --- New column G "X" marks (plain centered) for rows 9,10,11 ---
$ws.Range("G9").Value = "X"
$ws.Range("G9").HorizontalAlignment = -4108
$ws.Range("G10").Value = "X"
$ws.Range("G10").HorizontalAlignment = -4108
$ws.Range("G11").Value = "X"
$ws.Range("G11").HorizontalAlignment = -4108

# --- New column H (person assigned) values, centered ---
$ws.Range("H3").Value = "Julius"
$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("H8").Value = "Julius"
$ws.Range("H8").HorizontalAlignment = -4108
$ws.Range("H11").Value = "BWLer"
$ws.Range("H11").HorizontalAlignment = -4108

# --- New row content: B11 ---
$ws.Range("B11").Value = "Produktdatenbank mit Infos und Bildern füllen"

# --- Empty but centered/styled H placeholder cells (match column F pattern) ---
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H4").HorizontalAlignment = -4108
$ws.Range("H5").HorizontalAlignment = -4108
$ws.Range("H6").HorizontalAlignment = -4108
$ws.Range("H7").HorizontalAlignment = -4108
$ws.Range("H9").HorizontalAlignment = -4108
$ws.Range("H10").HorizontalAlignment = -4108
$ws.Range("H12").HorizontalAlignment = -4108

# --- Selection moved to B19 ---
$null = $ws.Range("B19").Select()
